# Updated symbol list on Thu Dec 29 13:40:22 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores its values as text that merely looks
# numeric (the source workbook uses inline strings everywhere). To keep
# these cells as text (instead of letting Excel silently convert a
# numeric-looking string into a real floating point number, which would
# also mangle trailing zeros such as "0.1330"), every Price update is
# written with a leading apostrophe (the classic "quote prefix" trick)
# so Excel stores it verbatim as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($addr, $text)
    $ws.Range($addr).Value = "'" + $text
}

# ---- straightforward price refreshes -------------------------------
Set-PriceText "D2"  "245.73"
Set-PriceText "D4"  "5.313"
Set-PriceText "D5"  "0.05737"
Set-PriceText "D6"  "6.507"
Set-PriceText "D7"  "3.135"
Set-PriceText "D8"  "0.8188"
Set-PriceText "D9"  "0.8693"
Set-PriceText "D10" "0.1376"
Set-PriceText "D11" "0.07009"
Set-PriceText "D13" "0.02923"
Set-PriceText "D14" "0.09398"
Set-PriceText "D15" "3.738"
Set-PriceText "D16" "0.001548"
Set-PriceText "D18" "0.0005974"
Set-PriceText "D19" "0.006195"
Set-PriceText "D20" "0.001242"
Set-PriceText "D21" "0.003860"
Set-PriceText "D22" "0.00008792"
Set-PriceText "D23" "3.541"
Set-PriceText "D24" "2.148"
Set-PriceText "D26" "0.1330"

# ---- row 28 (UpBots): price tweak + "Bestin24h" suffix removed ------
Set-PriceText "D28" "0.0003011"
$ws.Range("E28").Value = "27UpBotsUBXT"

Set-PriceText "D40" "0.03716"

# ---- rows 41-43: tokens rotate (Kick moves up, BKEXToken/CEJI shift
#      down by one row) and Kick gets a refreshed price + "Bestin24h" --
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-PriceText "D41" "0.006414"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-PriceText "D42" "0.1056"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-PriceText "D43" "0.002215"
$ws.Range("E43").Value = "42CEJICEJI"

# ---- remaining tail price refreshes ---------------------------------
Set-PriceText "D44" "0.008316"
Set-PriceText "D45" "0.00005212"
Set-PriceText "D47" "0.3596"
Set-PriceText "D48" "0.002253"
Set-PriceText "D49" "0.00002098"
Set-PriceText "D50" "0.0001998"
